$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: concept-type annotations. grupo / sexo / descripcion-ocupacion move
# from "dimension" to "measure"; direccion-provincial-nombre (column K) also
# becomes a measure instead of being tagged as the refArea dimension.
$ws.Range("G2").Value = "iaest-measure:grupo"
$ws.Range("H2").Value = "iaest-measure:sexo"
$ws.Range("I2").Value = "iaest-measure:descripcion-ocupacion"
$ws.Range("K2").Value = "iaest-measure:direccion-provincial-nombre"

# Row 3: medida/dim classification follows the same columns.
$ws.Range("G3").Value = "medida"
$ws.Range("H3").Value = "medida"
$ws.Range("I3").Value = "medida"
$ws.Range("K3").Value = "medida"

# Row 4: datatype follows suit - now plain ints instead of skos:Concept/URI.
$ws.Range("G4").Value = "xsd:int"
$ws.Range("H4").Value = "xsd:int"
$ws.Range("I4").Value = "xsd:int"
$ws.Range("K4").Value = "xsd:int"

# Row 5: mapping-file references for the now-retired dimension columns are
# no longer needed - clear the cells entirely (not just their value).
$ws.Range("G5").Clear()
$ws.Range("H5").Clear()
$ws.Range("I5").Clear()
